$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target stored widths: 26.42578125 / 19.42578125 / 4.85546875
# character-units, per the canonical OOXML <col> width attribute). This host's
# ColumnWidth setter quantizes to 1/6-character increments, so the inputs below
# are chosen to land on the closest achievable stored width to each target.
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
$ws.Columns.Item(3).ColumnWidth = 4

# Update mean mass flow rate values for 0.5s timestep run
$ws.Range("B2").Value = 0.02786632362762001
$ws.Range("B3").Value = 0.11523863538979752
$ws.Range("B4").Value = 0.097734540902479267
$ws.Range("B5").Value = 0.041750894314534381
$ws.Range("B6").Value = 0.21085593419090098
$ws.Range("B7").Value = 0.19900953139786906
$ws.Range("B8").Value = 0.069910496810435388
$ws.Range("B9").Value = 0.44772842357336207
